{"js": "// Load all tables in the document body.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The second table (\"Nom de balise / Champ correspondant / Format / ...\")\n// holds both rows that need changing.\nconst table = tables.items[1];\ntable.rows.load(\"items\");\nawait context.sync();\n\n// Load the cells for every row up front so we can inspect their text.\nfor (let i = 0; i < table.rows.items.length; i++) {\n  table.rows.items[i].cells.load(\"items\");\n}\nawait context.sync();\n\n// Load each cell's body text so we can find the target rows reliably\n// (by content) instead of relying purely on a fixed row index.\nconst rows = table.rows.items;\nfor (let i = 0; i < rows.length; i++) {\n  const cells = rows[i].cells.items;\n  for (let j = 0; j < cells.length; j++) {\n    cells[j].body.load(\"text\");\n  }\n}\nawait context.sync();\n\nlet centerNameRow = null;\nlet plateRow = null;\nfor (let i = 0; i < rows.length; i++) {\n  const cells = rows[i].cells.items;\n  const firstCellText = cells[0].body.text.trim();\n  if (firstCellText === \"centerName\") {\n    centerNameRow = rows[i];\n  } else if (firstCellText === \"plate\") {\n    plateRow = rows[i];\n  }\n}\n\n// 1) \"ID Centre d\u2019affectation\" -> \"Nom du centre d\u2019affectation\"\nif (centerNameRow) {\n  const cell = centerNameRow.cells.items[1];\n  cell.body.clear();\n  cell.body.insertText(\"Nom du centre d\\u2019affectation\", Word.InsertLocation.start);\n}\n\n// 2) Remove the whole \"plate\" / \"Immatriculation\" row.\nif (plateRow) {\n  plateRow.delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The table with \"Nom de balise / Champ correspondant / Format / Cardinalite\n# / Description / Exemple\" columns is the 2nd table in the document.\n$t = $d.Tables.Item(2)\n\n$centerNameRow = $null\n$plateRow = $null\nforeach ($row in $t.Rows) {\n    $firstCellText = $row.Cells.Item(1).Range.Text.TrimEnd([char]13, [char]7)\n    if ($firstCellText -eq \"centerName\") {\n        $centerNameRow = $row\n    } elseif ($firstCellText -eq \"plate\") {\n        $plateRow = $row\n    }\n}\n\n# 1) \"ID Centre d\u2019affectation\" -> \"Nom du centre d\u2019affectation\"\nif ($centerNameRow -ne $null) {\n    $centerNameRow.Cells.Item(2).Range.Text = \"Nom du centre d\u2019affectation\"\n}\n\n# 2) Remove the whole \"plate\" / \"Immatriculation\" row.\nif ($plateRow -ne $null) {\n    $plateRow.Delete()\n}\n"}
